# v1.4.0 sample data update: refresh clinical-event timeline on the
# "臨床イベント" (Clinical Events) sheet with the new, more granular
# symptom-level entries (thirst, polyuria, vomiting, abdominal pain, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("臨床イベント")
$ws.Activate()

# --- Rewrite existing rows 2-7 in place -----------------------------------
$data = @(
    @(2,  "2025-02-01", "口渇",     "著明な口渇、水分を頻回に要求"),
    @(3,  "2025-02-01", "多尿",     "頻尿、夜間尿増加の訴え"),
    @(4,  "2025-02-01", "嘔吐",     "頻回の嘔吐、経口摂取不可"),
    @(5,  "2025-02-01", "意識障害", "GCS E3V4M5、傾眠傾向"),
    @(6,  "2025-02-01", "脱水",     "皮膚ツルゴール低下、口唇乾燥著明"),
    @(7,  "2025-02-01", "頻呼吸",   "Kussmaul呼吸、RR 32/分"),
    @(8,  "2025-02-01", "腹痛",     "心窩部痛の訴え"),
    @(9,  "2025-02-01", "高血糖",   "血糖580mg/dL"),
    @(10, "2025-02-02", "意識障害", "GCS E4V5M6、意識清明に改善"),
    @(11, "2025-02-02", "高血糖",   "血糖165mg/dL、改善傾向"),
    @(12, "2025-02-02", "口渇",     "軽度改善"),
    @(13, "2025-02-03", "高血糖",   "血糖125mg/dL、正常化傾向")
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
}
